# Apply the recorded edits to the workbook:
#  1. Update the stored password hash on the Users sheet.
#  2. Strip every sheet except "Users" down to just the "id" header
#     (column A, row 1) — all other header columns and data rows removed.
#  3. Move the "Directory" sheet to the end of the tab order (after "Vessels").

$wb = $excel.ActiveWorkbook

# --- 1. Update Users!B2 (password_hash) -----------------------------------
$users = $wb.Worksheets.Item("Users")
$users.Range("B2").Value = "240be518fabd2724ddb6f04eeb1da5967448d7e831c08c8fa822809f74c720a9"

# --- 2. Trim every other sheet down to only the "id" header ---------------
$sheetsToTrim = @("Directory", "Requisitions", "Landings", "Categories", "Vessels")
foreach ($name in $sheetsToTrim) {
    $ws = $wb.Worksheets.Item($name)
    $ur = $ws.UsedRange
    $lastRow = $ur.Rows.Count
    $lastCol = $ur.Columns.Count

    if ($lastCol -gt 1) {
        $ws.Range($ws.Cells.Item(1, 2), $ws.Cells.Item($lastRow, $lastCol)).Clear()
    }
    if ($lastRow -gt 1) {
        $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 1)).Clear()
    }
}

# --- 3. Move "Directory" to the end of the tab order -----------------------
$directory = $wb.Worksheets.Item("Directory")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$directory.Move([System.Reflection.Missing]::Value, $lastSheet)
